# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 279
$wsOff.Range("C2").Value = 192
$wsOff.Range("D2").Value = 36

$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 324
$wsDef.Range("C2").Value = 206
$wsDef.Range("D2").Value = 88
$wsDef.Range("E2").Value = 43
$wsDef.Range("F2").Value = 8
